$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 69
$ws_ALC.Range("H69").Value = 4000
$ws_ALC.Range("I69").Value = 4000
$ws_ALC.Range("K69").Value = 12000
$ws_ALC.Range("M69").Value = -11126

# ALC row 70
$ws_ALC.Range("H70").Value = 0
$ws_ALC.Range("I70").Value = 0
$ws_ALC.Range("J70").Value = 0
$ws_ALC.Range("K70").Value = 0
$ws_ALC.Range("L70").Value = 0
$ws_ALC.Range("M70").ClearContents()
$ws_ALC.Range("N70").ClearContents()

# ALC row 72
$ws_ALC.Range("H72").Value = 4000
$ws_ALC.Range("I72").Value = 4000
$ws_ALC.Range("K72").Value = 36000
$ws_ALC.Range("M72").Value = -31632

# ALC row 73
$ws_ALC.Range("H73").Value = 0
$ws_ALC.Range("I73").Value = 0
$ws_ALC.Range("J73").Value = 0
$ws_ALC.Range("K73").Value = 0
$ws_ALC.Range("L73").Value = 0
$ws_ALC.Range("M73").ClearContents()
$ws_ALC.Range("N73").ClearContents()

# ALC row 74
$ws_ALC.Range("H74").Value = 1000
$ws_ALC.Range("I74").Value = 1000
$ws_ALC.Range("K74").Value = 1000
$ws_ALC.Range("M74").Value = -64

# ALC row 77
$ws_ALC.Range("H77").Value = 1000
$ws_ALC.Range("I77").Value = 1000
$ws_ALC.Range("K77").Value = 5000
$ws_ALC.Range("M77").Value = -320

# ALC row 80
$ws_ALC.Range("H80").Value = 3166.6667
$ws_ALC.Range("I80").Value = 5000
$ws_ALC.Range("J80").Value = 2250
$ws_ALC.Range("K80").Value = 15000
$ws_ALC.Range("L80").Value = 6750
$ws_ALC.Range("M80").Value = -14002
$ws_ALC.Range("N80").Value = -8746

# ALC row 83
$ws_ALC.Range("H83").Value = 3166.6667
$ws_ALC.Range("I83").Value = 5000
$ws_ALC.Range("J83").Value = 2250
$ws_ALC.Range("K83").Value = 45000
$ws_ALC.Range("L83").Value = 20250
$ws_ALC.Range("M83").Value = -40008
$ws_ALC.Range("N83").Value = -30234

# ALC row 93
$ws_ALC.Range("H93").Value = 0
$ws_ALC.Range("J93").Value = 0
$ws_ALC.Range("L93").Value = 0
$ws_ALC.Range("N93").ClearContents()

# ALC row 100
$ws_ALC.Range("H100").Value = 1269.6
$ws_ALC.Range("I100").Value = 1233
$ws_ALC.Range("K100").Value = 1233
$ws_ALC.Range("M100").Value = -692

# ARM row 38
$ws_ARM.Range("H38").Value = 1110773.1
$ws_ARM.Range("I38").Value = 5791.6
$ws_ARM.Range("K38").Value = 5791.6
$ws_ARM.Range("M38").Value = -5324.6

# ARM row 61
$ws_ARM.Range("H61").Value = 1716.6666
$ws_ARM.Range("I61").Value = 1716.6666
$ws_ARM.Range("K61").Value = 1716.6666
$ws_ARM.Range("M61").Value = -1504.6666

# ARM row 74
$ws_ARM.Range("H74").Value = 3250.5
$ws_ARM.Range("I74").Value = 2945
$ws_ARM.Range("J74").Value = 6000
$ws_ARM.Range("K74").Value = 2945
$ws_ARM.Range("L74").Value = 6000
$ws_ARM.Range("M74").Value = -2071
$ws_ARM.Range("N74").Value = -7748

# ARM row 77
$ws_ARM.Range("H77").Value = 3250.5
$ws_ARM.Range("I77").Value = 2945
$ws_ARM.Range("J77").Value = 6000
$ws_ARM.Range("K77").Value = 14725
$ws_ARM.Range("L77").Value = 30000
$ws_ARM.Range("M77").Value = -10357
$ws_ARM.Range("N77").Value = -38736

# ARM row 110
$ws_ARM.Range("H110").Value = 0
$ws_ARM.Range("I110").Value = 0
$ws_ARM.Range("J110").Value = 0
$ws_ARM.Range("K110").Value = 0
$ws_ARM.Range("L110").Value = 0
$ws_ARM.Range("M110").ClearContents()
$ws_ARM.Range("N110").ClearContents()

# ARM row 136
$ws_ARM.Range("H136").Value = 1716.6666
$ws_ARM.Range("I136").Value = 1716.6666
$ws_ARM.Range("K136").Value = 5149.9998
$ws_ARM.Range("M136").Value = -2599.9998

# BSM row 109
$ws_BSM.Range("H109").Value = 0
$ws_BSM.Range("J109").Value = 0
$ws_BSM.Range("L109").Value = 0
$ws_BSM.Range("N109").ClearContents()

# CRP row 132
$ws_CRP.Range("H132").Value = 2635.5
$ws_CRP.Range("I132").Value = 2453.5
$ws_CRP.Range("K132").Value = 7360.5
$ws_CRP.Range("M132").Value = -4830.5

# CUL row 38
$ws_CUL.Range("H38").Value = 1434.1666
$ws_CUL.Range("J38").Value = 2148.75
$ws_CUL.Range("L38").Value = 6446.25
$ws_CUL.Range("N38").Value = -7140.25

# CUL row 86
$ws_CUL.Range("H86").Value = 666.3333
$ws_CUL.Range("I86").Value = 249
$ws_CUL.Range("K86").Value = 747
$ws_CUL.Range("M86").Value = 439

# CUL row 89
$ws_CUL.Range("H89").Value = 666.3333
$ws_CUL.Range("I89").Value = 249
$ws_CUL.Range("K89").Value = 2241
$ws_CUL.Range("M89").Value = 3687

# GSM row 11
$ws_GSM.Range("H11").Value = 28917042
$ws_GSM.Range("J11").Value = 0
$ws_GSM.Range("L11").Value = 0
$ws_GSM.Range("N11").ClearContents()

# GSM row 122
$ws_GSM.Range("H122").Value = 2916.1667
$ws_GSM.Range("I122").Value = 1499.4
$ws_GSM.Range("K122").Value = 4498.200000000001
$ws_GSM.Range("M122").Value = -2048.200000000001

# GSM row 132
$ws_GSM.Range("H132").Value = 2644.5
$ws_GSM.Range("I132").Value = 1820
$ws_GSM.Range("J132").Value = 2997.8572
$ws_GSM.Range("K132").Value = 5460
$ws_GSM.Range("L132").Value = 8993.571599999999
$ws_GSM.Range("M132").Value = -2930
$ws_GSM.Range("N132").Value = -14053.5716

# LTW row 32
$ws_LTW.Range("H32").Value = 5422.6665
$ws_LTW.Range("I32").Value = 1507.2
$ws_LTW.Range("K32").Value = 1507.2
$ws_LTW.Range("M32").Value = -1190.2

# LTW row 40
$ws_LTW.Range("H40").Value = 633093.6
$ws_LTW.Range("I40").Value = 5001
$ws_LTW.Range("J40").Value = 1009949.2
$ws_LTW.Range("K40").Value = 5001
$ws_LTW.Range("L40").Value = 1009949.2
$ws_LTW.Range("M40").Value = -4865
$ws_LTW.Range("N40").Value = -1010221.2

# LTW row 46
$ws_LTW.Range("H46").Value = 2797
$ws_LTW.Range("I46").Value = 996.25
$ws_LTW.Range("J46").Value = 10000
$ws_LTW.Range("K46").Value = 996.25
$ws_LTW.Range("L46").Value = 10000
$ws_LTW.Range("M46").Value = -808.25
$ws_LTW.Range("N46").Value = -10376

# LTW row 68
$ws_LTW.Range("H68").Value = 0
$ws_LTW.Range("I68").Value = 0
$ws_LTW.Range("J68").Value = 0
$ws_LTW.Range("K68").Value = 0
$ws_LTW.Range("L68").Value = 0
$ws_LTW.Range("M68").ClearContents()
$ws_LTW.Range("N68").ClearContents()

# LTW row 71
$ws_LTW.Range("H71").Value = 0
$ws_LTW.Range("I71").Value = 0
$ws_LTW.Range("J71").Value = 0
$ws_LTW.Range("K71").Value = 0
$ws_LTW.Range("L71").Value = 0
$ws_LTW.Range("M71").ClearContents()
$ws_LTW.Range("N71").ClearContents()

# LTW row 82
$ws_LTW.Range("H82").Value = 1299.909
$ws_LTW.Range("I82").Value = 1322.1111
$ws_LTW.Range("J82").Value = 1200
$ws_LTW.Range("K82").Value = 1322.1111
$ws_LTW.Range("L82").Value = 1200
$ws_LTW.Range("M82").Value = -961.1111000000001
$ws_LTW.Range("N82").Value = -1922

# LTW row 85
$ws_LTW.Range("H85").Value = 1299.909
$ws_LTW.Range("I85").Value = 1322.1111
$ws_LTW.Range("J85").Value = 1200
$ws_LTW.Range("K85").Value = 1322.1111
$ws_LTW.Range("L85").Value = 1200
$ws_LTW.Range("M85").Value = -74.11110000000008
$ws_LTW.Range("N85").Value = -3696

# LTW row 132
$ws_LTW.Range("H132").Value = 1427.3684
$ws_LTW.Range("I132").Value = 1364.1765
$ws_LTW.Range("K132").Value = 4092.5295
$ws_LTW.Range("M132").Value = -1562.5295

# WVR row 23
$ws_WVR.Range("H23").Value = 1917.7142
$ws_WVR.Range("I23").Value = 304.8
$ws_WVR.Range("K23").Value = 304.8
$ws_WVR.Range("M23").Value = -75.80000000000001

# WVR row 62
$ws_WVR.Range("H62").Value = 0
$ws_WVR.Range("I62").Value = 0
$ws_WVR.Range("K62").Value = 0
$ws_WVR.Range("M62").ClearContents()

# WVR row 65
$ws_WVR.Range("H65").Value = 0
$ws_WVR.Range("I65").Value = 0
$ws_WVR.Range("K65").Value = 0
$ws_WVR.Range("M65").ClearContents()

# WVR row 100
$ws_WVR.Range("H100").Value = 642.3333
$ws_WVR.Range("I100").Value = 642.3333
$ws_WVR.Range("J100").Value = 0
$ws_WVR.Range("K100").Value = 1284.6666
$ws_WVR.Range("L100").Value = 0
$ws_WVR.Range("M100").Value = -743.6666
$ws_WVR.Range("N100").ClearContents()

# WVR row 132
$ws_WVR.Range("H132").Value = 1096.6364
$ws_WVR.Range("I132").Value = 1196.3
$ws_WVR.Range("J132").Value = 100
$ws_WVR.Range("K132").Value = 3588.9
$ws_WVR.Range("L132").Value = 300
$ws_WVR.Range("M132").Value = -1058.9
$ws_WVR.Range("N132").Value = -5360
